# Delete the row for "「アルマジロ」" (Excel row 689) from Sheet1.
# This shifts every row below it up by one, which matches the target diff
# (dimension shrinks from A1:C831 to A1:C830, and all rows from 690..831
# become 689..830).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(689).Delete()
